$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.951.01'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +1.79%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.505.93'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '602.25'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +2.38%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '173.26'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +2.62%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.45%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.499.78'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.51%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.21'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +7.34%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.581'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '46.41'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.067.20'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '8.28'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '606.22'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -1.68%  '
$ws.Range('B18').NumberFormat = "@"
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').NumberFormat = "@"
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.507.47'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('B19').NumberFormat = "@"
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').NumberFormat = "@"
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '70.021.88'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.64%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.91%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.12'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.79%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.871'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.37%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.08'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -18.15%  '
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -2.39%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '95.56'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -1.95%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.57'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -1.38%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '33.87'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +2.80%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.97'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -1.90%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '710.73'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +23.11%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.00'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -2.85%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '8.09'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -3.98%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.92'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +1.08%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -2.63%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0997'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -1.40%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.55'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '10.68'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -0.17%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0471'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +8.24%  '
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.18%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '56.41'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +4.30%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.327.99'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -2.44%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.312'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -3.52%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +4.11%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '32.23'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0₃0689'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.70%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.55'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.62%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +0.95%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '133.29'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.87%  '
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.00%  '

Write-Host "Applied 89 cell updates"
